$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has 3 item rows (7,8,9), a totals row (10) and a
# footer row (11). Two new sold items arrived, so we need to insert two
# more item rows (pushing the totals/footer rows down to 12/13), fill them
# in with the same look & feel as the existing item rows, update the
# running total and refresh the "printed at" timestamp in the footer.
# ---------------------------------------------------------------------------

# 1) Push the totals row (10) and footer row (11) down by two rows, making
#    room for the two new item rows at 10 and 11.
$ws.Rows("10:11").Insert()

# 2) Clone the formatting of an existing item row onto the two freshly
#    inserted (currently blank) rows so borders/fonts/fills/number formats
#    match the other item rows exactly.
$ws.Range("A7:Q7").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Item #4 -> row 10
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "بلاستر مترسيلك 2 سم"
$ws.Range("H10").Value = "27:0"

$fmt = $ws.Range("L10").NumberFormat
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "0"
$ws.Range("L10").NumberFormat = $fmt

$ws.Range("N10").Value = "15.00"

$fmt = $ws.Range("P10").NumberFormat
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "15.0000"
$ws.Range("P10").NumberFormat = $fmt

$ws.Range("Q10").Value = "1:0"

# 4) Item #5 -> row 11
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "كالونا "
$ws.Range("H11").Value = "0:0"

$fmt = $ws.Range("L11").NumberFormat
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = "0"
$ws.Range("L11").NumberFormat = $fmt

$ws.Range("N11").Value = "15.00"

$fmt = $ws.Range("P11").NumberFormat
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "15.0000"
$ws.Range("P11").NumberFormat = $fmt

$ws.Range("Q11").Value = "1:0"

# 5) Re-create the merged cell layout used by the other item rows for the
#    two new rows (PasteSpecial of formats only does not carry merges over).
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# 6) The totals row (now row 12) needs to reflect the two new sales
#    (260 + 15 + 15 = 290).
$ws.Range("P12").Value = 290

# 7) The footer (now row 13) keeps its page-number / credit text, only the
#    "printed at" timestamp advances to the new export time.
$ws.Range("A13").Value = "Sunday, 7 September, 2025 9:44 AM"

Write-Output "Inserted 2 new sale rows and refreshed totals/timestamp."
